$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10: Pause (typed first, per shared-string insertion order)
$ws.Range("A10").Value = "Pause"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0

# New rows 5-6: LeftFast / LeftSlow
$ws.Range("A5").Value = "LeftFast"
$ws.Range("B5").Value = -2
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1

$ws.Range("A6").Value = "LeftSlow"
$ws.Range("B6").Value = -2
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 3

# New rows 8-9: RightFast / RightSlow
$ws.Range("A8").Value = "RightFast"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 1

$ws.Range("A9").Value = "RightSlow"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 3

# Header rename: DurationTime -> Speed, 완료시간(이동시간) -> 이동속도
$ws.Range("D1").Value = "Speed"
$ws.Range("D2").Value = "이동속도"

# Existing rows kept/updated values
$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "DestX"
$ws.Range("C1").Value = "DestY"

$ws.Range("A3").Value = "string"
$ws.Range("B3").Value = "float"
$ws.Range("C3").Value = "float"
$ws.Range("D3").Value = "int"

$ws.Range("A4").Value = "Left"
$ws.Range("B4").Value = -2
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 10

$ws.Range("A7").Value = "Right"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 2

$ws.Range("D4").Select()
